$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now reflects former row 4 data
$ws.Range("D2").Value = 44701
$ws.Range("K2").Value = 19000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 19500
$ws.Range("P2").Value = 650

# Row 3 now reflects former row 2 data
$ws.Range("D3").Value = 44438
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("O3").Value = 'Provincia del Elquí'
$ws.Range("P3").Value = 383

# Row 4 now reflects former row 8 data
$ws.Range("D4").Value = 44427
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12500
$ws.Range("N4").Value = '$/caja 40 unidades'
$ws.Range("O4").Value = 'Provincia de Limarí'
$ws.Range("P4").Value = 312
$ws.Range("Q4").Value = 40

# Row 8 now reflects former row 3 data
$ws.Range("D8").Value = 44498
$ws.Range("H8").Value = 'Española'
$ws.Range("K8").Value = 8500
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8750
$ws.Range("N8").Value = '$/caja 30 unidades'
$ws.Range("P8").Value = 292
$ws.Range("Q8").Value = 30

# Row 9 now reflects former row 10 data
$ws.Range("D9").Value = 44420
$ws.Range("H9").Value = 'Madrigal'
$ws.Range("J9").Value = 800
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 14500
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("P9").Value = 362
$ws.Range("Q9").Value = 40

# Row 10 now reflects former row 11 data
$ws.Range("J10").Value = 700
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13500
$ws.Range("O10").Value = 'Provincia del Elquí'
$ws.Range("P10").Value = 338

# Row 11 now reflects former row 9 data
$ws.Range("D11").Value = 44687
$ws.Range("H11").Value = 'Española'
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 19000
$ws.Range("M11").Value = 18500
$ws.Range("N11").Value = '$/caja 30 unidades'
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 617
$ws.Range("Q11").Value = 30
